$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @(
    @{ Row = 82; I = "b"; J = "Acknowledge (Backchannel)" }
    @{ Row = 90; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 96; I = "sv"; J = "Statement-opinion" }
    @{ Row = 108; I = "ba"; J = "Appreciation" }
    @{ Row = 111; I = "b"; J = "Acknowledge (Backchannel)" }
    @{ Row = 113; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 115; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 120; I = "b"; J = "Acknowledge (Backchannel)" }
    @{ Row = 159; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 160; I = "aa"; J = "Agree/Accept" }
    @{ Row = 171; I = "ba"; J = "Appreciation" }
    @{ Row = 175; I = "ba"; J = "Appreciation" }
    @{ Row = 179; I = "ba"; J = "Appreciation" }
    @{ Row = 192; I = "sv"; J = "Statement-opinion" }
    @{ Row = 219; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 221; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 222; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 230; I = "sv"; J = "Statement-opinion" }
    @{ Row = 238; I = "b"; J = "Acknowledge (Backchannel)" }
    @{ Row = 252; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 253; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 259; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 271; I = "sv"; J = "Statement-opinion" }
    @{ Row = 293; I = "sv"; J = "Statement-opinion" }
    @{ Row = 310; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 313; I = "sv"; J = "Statement-opinion" }
    @{ Row = 317; I = "b"; J = "Acknowledge (Backchannel)" }
    @{ Row = 323; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 341; I = "sv"; J = "Statement-opinion" }
    @{ Row = 363; I = "%"; J = "Uninterpretable" }
    @{ Row = 372; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 373; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 409; I = "aa"; J = "Agree/Accept" }
    @{ Row = 420; I = "aa"; J = "Agree/Accept" }
    @{ Row = 421; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 426; I = "b"; J = "Acknowledge (Backchannel)" }
    @{ Row = 436; I = "sv"; J = "Statement-opinion" }
    @{ Row = 437; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 439; I = "sv"; J = "Statement-opinion" }
    @{ Row = 440; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 459; I = "sv"; J = "Statement-opinion" }
    @{ Row = 478; I = "sv"; J = "Statement-opinion" }
    @{ Row = 481; I = "ba"; J = "Appreciation" }
    @{ Row = 497; I = "b"; J = "Acknowledge (Backchannel)" }
    @{ Row = 504; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 510; I = "ba"; J = "Appreciation" }
    @{ Row = 552; I = "aa"; J = "Agree/Accept" }
    @{ Row = 559; I = "aa"; J = "Agree/Accept" }
    @{ Row = 560; I = "sv"; J = "Statement-opinion" }
    @{ Row = 563; I = "sv"; J = "Statement-opinion" }
    @{ Row = 579; I = "sd"; J = "Statement-non-opinion" }
)

foreach ($u in $updates) {
    $ws.Range("I" + $u.Row).Value = $u.I
    $ws.Range("J" + $u.Row).Value = $u.J
}
